$wb = $excel.ActiveWorkbook

# --- 1. Create the new "RC_P10" sheet as a copy of "RC_P9" (placed right after it) ---
$src = $wb.Worksheets.Item("RC_P9")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "RC_P10"

# --- 2. Update the header block (B1:B6) for the new cohort / collection ---
$ws.Range("B1").Value = '2018-10-14'
$ws.Range("B2").Value = 'cohort day 9 (2018-09-29) apple + haw cohort; cohort day 10 (2018-09-30) apple + haw host; cohort day 11 (2018-10-01) apple + haw; cohort day 12 ( 2018-10-02) apple + haw '
$ws.Range("B4").Value = '10'
$ws.Range("B5").Value = 'Chao''s rack'
$ws.Range("B6").Value = 'nb#005, pg 125, 129 , 114, 119'

# --- 3. Fill in the box-position grid (rows 9-13) with the new vial labels ---
$ws.Range("C9").Value = 'H-9-1'
$ws.Range("D9").Value = 'H-9-2'
$ws.Range("E9").Value = 'H-9-3'
$ws.Range("F9").Value = 'H-9-4'
$ws.Range("G9").Value = 'H-9-5'
$ws.Range("H9").Value = 'H-9-6'
$ws.Range("I9").Value = 'H-9-7'
$ws.Range("J9").Value = 'H-9-8'
$ws.Range("K9").Value = 'H-9-9'

$ws.Range("C10").Value = 'H-9-10'
$ws.Range("D10").Value = 'H-10-1'
$ws.Range("E10").Value = 'H-10-2'
$ws.Range("F10").Value = 'H-10-3'
$ws.Range("G10").Value = 'H-10-4'
$ws.Range("H10").Value = 'H-10-5'
$ws.Range("I10").Value = 'H-10-6'
$ws.Range("J10").Value = 'H-10-7'
$ws.Range("K10").Value = 'H-10-8'

$ws.Range("C11").Value = 'H-10-9'
$ws.Range("D11").Value = 'H-10-10'
$ws.Range("E11").Value = 'A3-10-1'
$ws.Range("F11").Value = 'A3-10-2'
$ws.Range("G11").Value = 'A3-10-3'
$ws.Range("H11").Value = 'A3-10-4'
$ws.Range("I11").Value = 'H-11-1'
$ws.Range("J11").Value = 'H-11-2'
$ws.Range("K11").Value = 'H-11-3'

$ws.Range("C12").Value = 'H-11-4'
$ws.Range("D12").Value = 'H-11-5'
$ws.Range("E12").Value = 'A3-11-1'
$ws.Range("F12").Value = 'A3-11-2'
$ws.Range("G12").Value = 'A3-11-3'
$ws.Range("H12").Value = 'A3-12-1'
$ws.Range("I12").Value = 'A3-12-2'
$ws.Range("J12").Value = 'A3-12-3'
$ws.Range("K12").Value = 'A3-12-4'

$ws.Range("C13").Value = 'H-12-1'
$ws.Range("D13").Value = 'H-12-2'
$ws.Range("E13").Value = 'H-12-3'
$ws.Range("F13").Value = 'H-12-4'

# --- 4. Fix up cell shading so the new rows match the color-banding pattern ---
$ws.Range("C9").Copy()
$ws.Range("D9:K9").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("K10").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)

# --- 5. Make RC_P10 the active tab, mirroring where the author left off ---
$ws.Activate()
$ws.Range("G13").Select()
